$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2 "I" -> "WOMAN1", D2 stays the same text (reorders shared strings naturally)
$ws.Range("B2").Value = "WOMAN1"

# New row 3: cognate set entry TG100 / PERSON1 / /am/ (description) {anysource}
$ws.Range("A3").Value = $ws.Range("A2").Value2
$ws.Range("B3").Value = "PERSON1"
$ws.Range("D3").Value = "/am/ (description) {anysource}"

# Update the comment on A2 to describe it as a cognate set comment
$ws.Range("A2").Comment.Text("Cognate set comment")

# Move the active selection, matching the saved view state
$ws.Range("D14").Select()
